# [pr4-2896] Harmonize Configuration schema with the current information model.
#
# The "Compensation Good" lookup values on the PCNCode sheet are renamed to
# their new (capitalised) spellings, and the PCNCode sheet becomes the
# workbook's active tab (previously CustomsUnion was active).

$wb = $excel.ActiveWorkbook

# --- PCNCode sheet: rename the "Compensation Good" column values ----------
$pcn = $wb.Worksheets.Item("PCNCode")

$pcn.Range("C2").Value = "Papierosy"
$pcn.Range("C3").Value = "PyłTytoiowy"
$pcn.Range("C4").Value = "Kartony"
$pcn.Range("C5").Value = "Krajanka"

# --- Make PCNCode the active sheet / tab, with the new selection ----------
# (this also clears tabSelected on the previously-active CustomsUnion sheet)
$pcn.Activate() | Out-Null

$view = $excel.ActiveWindow
$view.ScrollRow = 4
$view.ScrollColumn = 1

$pcn.Range("C9").Select() | Out-Null
